# Natmi following Dr Hou advice
# Recompute Ligand/Receptor-expressing cell counts (now 3 instead of 1 per
# cluster) and all downstream derived expression/specificity metrics for
# rows 2-17 of the Calr-Scarf1 LR-pairs sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    ,(2, 5, 3)
    ,(2, 7, 81.03634266666667)
    ,(2, 8, 243.109028)
    ,(2, 9, 0.1632931649012984)
    ,(2, 10, 0.1632931649012984)
    ,(2, 11, 3)
    ,(2, 13, 9.575242333333334)
    ,(2, 14, 28.725727)
    ,(2, 15, 0.6572152168068958)
    ,(2, 16, 0.6572152168068958)
    ,(2, 17, 775.9426188403729)
    ,(2, 18, 6983.483569563356)
    ,(2, 19, 0.107318752773691)
    ,(2, 20, 0.107318752773691)
    ,(3, 5, 3)
    ,(3, 7, 81.03634266666667)
    ,(3, 8, 243.109028)
    ,(3, 9, 0.1632931649012984)
    ,(3, 10, 0.1632931649012984)
    ,(3, 11, 3)
    ,(3, 13, 0.34342)
    ,(3, 14, 1.03026)
    ,(3, 15, 0.0235712937488918)
    ,(3, 16, 0.0235712937488918)
    ,(3, 17, 27.82950079858667)
    ,(3, 18, 250.46550718728)
    ,(3, 19, 0.003849031157074732)
    ,(3, 20, 0.003849031157074732)
    ,(4, 5, 3)
    ,(4, 7, 81.03634266666667)
    ,(4, 8, 243.109028)
    ,(4, 9, 0.1632931649012984)
    ,(4, 10, 0.1632931649012984)
    ,(4, 11, 3)
    ,(4, 13, 2.978372666666667)
    ,(4, 14, 8.935117999999999)
    ,(4, 15, 0.2044263497165867)
    ,(4, 16, 0.2044263497165867)
    ,(4, 17, 241.3564280050338)
    ,(4, 18, 2172.207852045304)
    ,(4, 19, 0.03338142563444108)
    ,(4, 20, 0.03338142563444108)
    ,(5, 5, 3)
    ,(5, 7, 81.03634266666667)
    ,(5, 8, 243.109028)
    ,(5, 9, 0.1632931649012984)
    ,(5, 10, 0.1632931649012984)
    ,(5, 11, 3)
    ,(5, 13, 1.672381666666666)
    ,(5, 14, 5.017144999999999)
    ,(5, 15, 0.1147871397276258)
    ,(5, 16, 0.1147871397276258)
    ,(5, 17, 135.5236938094511)
    ,(5, 18, 1219.71324428506)
    ,(5, 19, 0.01874395533609157)
    ,(5, 20, 0.01874395533609157)
    ,(6, 5, 3)
    ,(6, 7, 178.365814)
    ,(6, 8, 535.097442)
    ,(6, 9, 0.3594179761796791)
    ,(6, 10, 0.3594179761796791)
    ,(6, 11, 3)
    ,(6, 13, 9.575242333333334)
    ,(6, 14, 28.725727)
    ,(6, 15, 0.6572152168068958)
    ,(6, 16, 0.6572152168068958)
    ,(6, 17, 1707.895893032259)
    ,(6, 18, 15371.06303729033)
    ,(6, 19, 0.2362149631392235)
    ,(6, 20, 0.2362149631392235)
    ,(7, 5, 3)
    ,(7, 7, 178.365814)
    ,(7, 8, 535.097442)
    ,(7, 9, 0.3594179761796791)
    ,(7, 10, 0.3594179761796791)
    ,(7, 11, 3)
    ,(7, 13, 0.34342)
    ,(7, 14, 1.03026)
    ,(7, 15, 0.0235712937488918)
    ,(7, 16, 0.0235712937488918)
    ,(7, 17, 61.25438784388)
    ,(7, 18, 551.2894905949199)
    ,(7, 19, 0.008471946695163411)
    ,(7, 20, 0.008471946695163412)
    ,(8, 5, 3)
    ,(8, 7, 178.365814)
    ,(8, 8, 535.097442)
    ,(8, 9, 0.3594179761796791)
    ,(8, 10, 0.3594179761796791)
    ,(8, 11, 3)
    ,(8, 13, 2.978372666666667)
    ,(8, 14, 8.935117999999999)
    ,(8, 15, 0.2044263497165867)
    ,(8, 16, 0.2044263497165867)
    ,(8, 17, 531.2398650853506)
    ,(8, 18, 4781.158785768155)
    ,(8, 19, 0.07347450489293489)
    ,(8, 20, 0.0734745048929349)
    ,(9, 5, 3)
    ,(9, 7, 178.365814)
    ,(9, 8, 535.097442)
    ,(9, 9, 0.3594179761796791)
    ,(9, 10, 0.3594179761796791)
    ,(9, 11, 3)
    ,(9, 13, 1.672381666666666)
    ,(9, 14, 5.017144999999999)
    ,(9, 15, 0.1147871397276258)
    ,(9, 16, 0.1147871397276258)
    ,(9, 17, 298.2957172936766)
    ,(9, 18, 2684.66145564309)
    ,(9, 19, 0.04125656145235729)
    ,(9, 20, 0.0412565614523573)
    ,(10, 5, 3)
    ,(10, 7, 146.2303303333333)
    ,(10, 8, 438.6909910000001)
    ,(10, 9, 0.2946630198121519)
    ,(10, 10, 0.2946630198121519)
    ,(10, 11, 3)
    ,(10, 13, 9.575242333333334)
    ,(10, 14, 28.725727)
    ,(10, 15, 0.6572152168068958)
    ,(10, 16, 0.6572152168068958)
    ,(10, 17, 1400.190849425051)
    ,(10, 18, 12601.71764482546)
    ,(10, 19, 0.1936570204508181)
    ,(10, 20, 0.1936570204508181)
    ,(11, 5, 3)
    ,(11, 7, 146.2303303333333)
    ,(11, 8, 438.6909910000001)
    ,(11, 9, 0.2946630198121519)
    ,(11, 10, 0.2946630198121519)
    ,(11, 11, 3)
    ,(11, 13, 0.34342)
    ,(11, 14, 1.03026)
    ,(11, 15, 0.0235712937488918)
    ,(11, 16, 0.0235712937488918)
    ,(11, 17, 50.21842004307334)
    ,(11, 18, 451.9657803876601)
    ,(11, 19, 0.006945588596927757)
    ,(11, 20, 0.006945588596927758)
    ,(12, 5, 3)
    ,(12, 7, 146.2303303333333)
    ,(12, 8, 438.6909910000001)
    ,(12, 9, 0.2946630198121519)
    ,(12, 10, 0.2946630198121519)
    ,(12, 11, 3)
    ,(12, 13, 2.978372666666667)
    ,(12, 14, 8.935117999999999)
    ,(12, 15, 0.2044263497165867)
    ,(12, 16, 0.2044263497165867)
    ,(12, 17, 435.5284189024376)
    ,(12, 18, 3919.755770121938)
    ,(12, 19, 0.06023688553666447)
    ,(12, 20, 0.06023688553666449)
    ,(13, 5, 3)
    ,(13, 7, 146.2303303333333)
    ,(13, 8, 438.6909910000001)
    ,(13, 9, 0.2946630198121519)
    ,(13, 10, 0.2946630198121519)
    ,(13, 11, 3)
    ,(13, 13, 1.672381666666666)
    ,(13, 14, 5.017144999999999)
    ,(13, 15, 0.1147871397276258)
    ,(13, 16, 0.1147871397276258)
    ,(13, 17, 244.5529235600772)
    ,(13, 18, 2200.976312040695)
    ,(13, 19, 0.03382352522774165)
    ,(13, 20, 0.03382352522774165)
    ,(14, 5, 3)
    ,(14, 7, 90.630432)
    ,(14, 8, 271.891296)
    ,(14, 9, 0.1826258391068707)
    ,(14, 10, 0.1826258391068707)
    ,(14, 11, 3)
    ,(14, 13, 9.575242333333334)
    ,(14, 14, 28.725727)
    ,(14, 15, 0.6572152168068958)
    ,(14, 16, 0.6572152168068958)
    ,(14, 17, 867.808349174688)
    ,(14, 18, 7810.275142572192)
    ,(14, 19, 0.1200244804431633)
    ,(14, 20, 0.1200244804431633)
    ,(15, 5, 3)
    ,(15, 7, 90.630432)
    ,(15, 8, 271.891296)
    ,(15, 9, 0.1826258391068707)
    ,(15, 10, 0.1826258391068707)
    ,(15, 11, 3)
    ,(15, 13, 0.34342)
    ,(15, 14, 1.03026)
    ,(15, 15, 0.0235712937488918)
    ,(15, 16, 0.0235712937488918)
    ,(15, 17, 31.12430295744)
    ,(15, 18, 280.11872661696)
    ,(15, 19, 0.004304727299725901)
    ,(15, 20, 0.004304727299725902)
    ,(16, 5, 3)
    ,(16, 7, 90.630432)
    ,(16, 8, 271.891296)
    ,(16, 9, 0.1826258391068707)
    ,(16, 10, 0.1826258391068707)
    ,(16, 11, 3)
    ,(16, 13, 2.978372666666667)
    ,(16, 14, 8.935117999999999)
    ,(16, 15, 0.2044263497165867)
    ,(16, 16, 0.2044263497165867)
    ,(16, 17, 269.931201436992)
    ,(16, 18, 2429.380812932928)
    ,(16, 19, 0.03733353365254623)
    ,(16, 20, 0.03733353365254625)
    ,(17, 5, 3)
    ,(17, 7, 90.630432)
    ,(17, 8, 271.891296)
    ,(17, 9, 0.1826258391068707)
    ,(17, 10, 0.1826258391068707)
    ,(17, 11, 3)
    ,(17, 13, 1.672381666666666)
    ,(17, 14, 5.017144999999999)
    ,(17, 15, 0.1147871397276258)
    ,(17, 16, 0.1147871397276258)
    ,(17, 17, 151.56867291888)
    ,(17, 18, 1364.11805626992)
    ,(17, 19, 0.02096309771143527)
    ,(17, 20, 0.02096309771143527)
)

foreach ($item in $changes) {
    $r = $item[0]
    $c = $item[1]
    $v = $item[2]
    $ws.Cells.Item($r, $c).Value = $v
}
